$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "spr_placeholder" sprite references with the real sprite
# names for the newly added Severe/Critical potions. Setting them in this
# order matches the order new shared-string entries get appended so the
# resulting sharedStrings.xml lines up with the target workbook.
$ws.Range("B8").Value = "spr_armor_potion_severe"
$ws.Range("B4").Value = "spr_health_potion_severe"
$ws.Range("B5").Value = "spr_health_potion_critical"
$ws.Range("B9").Value = "spr_armor_potion_critical"

# Update the selected cell to match the saved view state.
[void]$ws.Range("C15").Select()
